$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as plain text (inline/shared string),
# matching the source data which always stores Price/Volume as text cells,
# even when the text looks like a plain number (e.g. "1.005").
function Set-TextCell($range, $text) {
    $range.Formula = "'" + $text
    $range.Style = "Normal"
}

# Update Price (D) and Volume(1h) (E) columns for rows 2-49 based on refreshed crypto data
Set-TextCell $ws.Range('D2') '26.402.27'
Set-TextCell $ws.Range('E2') '  -7.31%  '
Set-TextCell $ws.Range('D3') '1.684.90'
Set-TextCell $ws.Range('E3') '  -5.94%  '
Set-TextCell $ws.Range('D4') '1.005'
Set-TextCell $ws.Range('E4') '  +0.14%  '
Set-TextCell $ws.Range('D5') '219.19'
Set-TextCell $ws.Range('E5') '  -5.18%  '
Set-TextCell $ws.Range('D6') '0.5068'
Set-TextCell $ws.Range('E6') '  -13.72%  '
Set-TextCell $ws.Range('D7') '1.005'
Set-TextCell $ws.Range('E7') '  +0.04%  '
Set-TextCell $ws.Range('D8') '0.2667'
Set-TextCell $ws.Range('E8') '  -3.52%  '
Set-TextCell $ws.Range('D9') '22.06'
Set-TextCell $ws.Range('E9') '  -5.22%  '
Set-TextCell $ws.Range('D10') '0.06305'
Set-TextCell $ws.Range('E10') '  -6.31%  '
Set-TextCell $ws.Range('D11') '0.07392'
Set-TextCell $ws.Range('E11') '  -2.26%  '
Set-TextCell $ws.Range('D12') '1.686.64'
Set-TextCell $ws.Range('E12') '  -6.07%  '
Set-TextCell $ws.Range('D13') '4.532'
Set-TextCell $ws.Range('E13') '  -5.22%  '
Set-TextCell $ws.Range('D14') '0.5785'
Set-TextCell $ws.Range('E14') '  -5.18%  '
Set-TextCell $ws.Range('D15') '1.918.77'
Set-TextCell $ws.Range('E15') '  -5.73%  '
Set-TextCell $ws.Range('D16') '0.000008644'
Set-TextCell $ws.Range('E16') '  -2.13%  '
Set-TextCell $ws.Range('D17') '65.17'
Set-TextCell $ws.Range('E17') '  -13.65%  '
Set-TextCell $ws.Range('D18') '26.481.47'
Set-TextCell $ws.Range('E18') '  -7.03%  '
Set-TextCell $ws.Range('D19') '5.003'
Set-TextCell $ws.Range('E19') '  -7.82%  '
Set-TextCell $ws.Range('E20') '  +0.05%  '
Set-TextCell $ws.Range('D21') '10.89'
Set-TextCell $ws.Range('E21') '  -4.65%  '
Set-TextCell $ws.Range('D22') '186.90'
Set-TextCell $ws.Range('E22') '  -10.00%  '
Set-TextCell $ws.Range('D23') '6.265'
Set-TextCell $ws.Range('E23') '  -7.90%  '
Set-TextCell $ws.Range('D24') '1.006'
Set-TextCell $ws.Range('E24') '  +0.08%  '
Set-TextCell $ws.Range('D25') '144.59'
Set-TextCell $ws.Range('E25') '  -5.27%  '
Set-TextCell $ws.Range('D26') '7.495'
Set-TextCell $ws.Range('E26') '  -5.95%  '
Set-TextCell $ws.Range('E27') '  -6.88%  '
Set-TextCell $ws.Range('D28') '15.90'
Set-TextCell $ws.Range('E28') '  -2.89%  '
Set-TextCell $ws.Range('D29') '1.347'
Set-TextCell $ws.Range('E29') '  -4.85%  '
Set-TextCell $ws.Range('D30') '0.05746'
Set-TextCell $ws.Range('E30') '  -5.81%  '
Set-TextCell $ws.Range('D31') '1.332'
Set-TextCell $ws.Range('E31') '  -6.09%  '
Set-TextCell $ws.Range('D32') '3.520'
Set-TextCell $ws.Range('E32') '  -6.81%  '
Set-TextCell $ws.Range('D33') '3.521'
Set-TextCell $ws.Range('E33') '  -6.33%  '
Set-TextCell $ws.Range('D34') '1.667'
Set-TextCell $ws.Range('E34') '  -2.99%  '
Set-TextCell $ws.Range('D35') '1.013'
Set-TextCell $ws.Range('E35') '  -3.16%  '
Set-TextCell $ws.Range('D36') '0.5969'
Set-TextCell $ws.Range('E36') '  -6.96%  '
Set-TextCell $ws.Range('D37') '2.356'
Set-TextCell $ws.Range('E37') '  -5.78%  '
Set-TextCell $ws.Range('D38') '2.681'
Set-TextCell $ws.Range('E38') '  -0.82%  '
Set-TextCell $ws.Range('D39') '1.099.81'
Set-TextCell $ws.Range('E39') '  -4.28%  '
Set-TextCell $ws.Range('D40') '0.01605'
Set-TextCell $ws.Range('E40') '  -4.35%  '
Set-TextCell $ws.Range('D41') '5.893'
Set-TextCell $ws.Range('E41') '  -6.47%  '
Set-TextCell $ws.Range('D42') '0.8604'
Set-TextCell $ws.Range('E42') '  -1.40%  '
Set-TextCell $ws.Range('D43') '1.003'
Set-TextCell $ws.Range('E43') '  -0.10%  '
Set-TextCell $ws.Range('D44') '99.81'
Set-TextCell $ws.Range('E44') '  -0.70%  '
Set-TextCell $ws.Range('D45') '1.845.39'
Set-TextCell $ws.Range('E45') '  -5.22%  '
Set-TextCell $ws.Range('D46') '0.00000000114'
Set-TextCell $ws.Range('E46') '  +3.84%  '
Set-TextCell $ws.Range('D47') '56.47'
Set-TextCell $ws.Range('E47') '  -6.09%  '
Set-TextCell $ws.Range('E48') '  +0.69%  '
Set-TextCell $ws.Range('D49') '8.031'
Set-TextCell $ws.Range('E49') '  -3.77%  '

# Rows 50-51: Mantle/Cronos swapped display order, each with refreshed price and volume figures
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Range('D50') '0.05222'
Set-TextCell $ws.Range('E50') '  -3.92%  '

$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell $ws.Range('D51') '0.4311'
Set-TextCell $ws.Range('E51') '  -3.66%  '
